$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($addr in @("D5", "D6", "D10", "D12", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D30", "D31", "D32", "D33", "D34", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.339.81"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").Value = "2.179.78"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "254.00"
$ws.Range("D6").Value = "0.610"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +2.04%  "
$ws.Range("D10").Value = "40.19"
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("D12").Value = "6.80"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").Value = "2.508.82"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").Value = "14.30"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "2.177.26"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").Value = "42.262.44"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D19").Value = "0.0000102"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "70.97"
$ws.Range("E20").Value = "  +1.63%  "
$ws.Range("D21").Value = "5.89"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").Value = "2.18"
$ws.Range("E22").Value = "  +6.30%  "
$ws.Range("D23").Value = "227.60"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").Value = "9.59"
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").Value = "10.58"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").Value = "3.32"
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("E29").Value = "  +2.86%  "
$ws.Range("D30").Value = "36.95"
$ws.Range("E30").Value = "  +12.58%  "
$ws.Range("D31").Value = "168.86"
$ws.Range("E31").Value = "  -0.67%  "
$ws.Range("D32").Value = "20.04"
$ws.Range("E32").Value = "  +1.40%  "
$ws.Range("D33").Value = "0.0809"
$ws.Range("E33").Value = "  +5.01%  "
$ws.Range("D34").Value = "5.15"
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("E36").Value = "  +2.60%  "
$ws.Range("D37").Value = "4.30"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "0.0329"
$ws.Range("E38").Value = "  +8.56%  "
$ws.Range("D39").Value = "12.04"
$ws.Range("E39").Value = "  +0.74%  "
$ws.Range("D40").Value = "2.07"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("E41").Value = "  +5.13%  "
$ws.Range("D42").Value = "5.20"
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("D43").Value = "57.85"
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("D44").Value = "102.25"
$ws.Range("E44").Value = "  +5.33%  "
$ws.Range("D45").Value = "0.473"
$ws.Range("E45").Value = "  +21.07%  "
$ws.Range("D46").Value = "8.31"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").Value = "0.0972"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("D48").Value = "2.42"
$ws.Range("E48").Value = "  +11.09%  "
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("E50").Value = "  +1.38%  "
$ws.Range("E51").Value = "  +1.23%  "
